# Auto-generated Excel COM-interop script
# Applies scheduled-runner price/profit updates to the Brynhildr_Profits workbook.
# For each affected sheet/row, refreshes currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ)
# and LeveProfit(NQ/HQ) columns (H-N) with refreshed market-board figures.

$wb = $excel.ActiveWorkbook

# ---- ALC sheet ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 241
$ws.Range("I2").Value = 54.666668
$ws.Range("K2").Value = 54.666668
$ws.Range("M2").Value = 58.333332
$ws.Range("H8").Value = 26.222221
$ws.Range("I8").Value = 26.222221
$ws.Range("K8").Value = 78.666663
$ws.Range("M8").Value = 60.333337
$ws.Range("H9").Value = 238.5
$ws.Range("I9").Value = 89.09999999999999
$ws.Range("J9").Value = 612
$ws.Range("K9").Value = 89.09999999999999
$ws.Range("L9").Value = 612
$ws.Range("M9").Value = 79.90000000000001
$ws.Range("N9").Value = -950
$ws.Range("H40").Value = 1477.2222
$ws.Range("I40").Value = 1474.375
$ws.Range("J40").Value = 1500
$ws.Range("K40").Value = 1474.375
$ws.Range("L40").Value = 1500
$ws.Range("M40").Value = -1299.375
$ws.Range("N40").Value = -1850
$ws.Range("H70").Value = 2460.375
$ws.Range("I70").Value = 2297
$ws.Range("J70").Value = 2558.4
$ws.Range("K70").Value = 6891
$ws.Range("L70").Value = 7675.200000000001
$ws.Range("M70").Value = -6621
$ws.Range("N70").Value = -8215.200000000001
$ws.Range("H73").Value = 2460.375
$ws.Range("I73").Value = 2297
$ws.Range("J73").Value = 2558.4
$ws.Range("K73").Value = 6891
$ws.Range("L73").Value = 7675.200000000001
$ws.Range("M73").Value = -5955
$ws.Range("N73").Value = -9547.200000000001
$ws.Range("H86").Value = 7316.6665
$ws.Range("I86").Value = 8500
$ws.Range("J86").Value = 4950
$ws.Range("K86").Value = 8500
$ws.Range("L86").Value = 4950
$ws.Range("M86").Value = -7377
$ws.Range("N86").Value = -7196
$ws.Range("H89").Value = 7316.6665
$ws.Range("I89").Value = 8500
$ws.Range("J89").Value = 4950
$ws.Range("K89").Value = 42500
$ws.Range("L89").Value = 24750
$ws.Range("M89").Value = -36884
$ws.Range("N89").Value = -35982
$ws.Range("H98").Value = 1963.5
$ws.Range("I98").Value = 2079.1428
$ws.Range("K98").Value = 2079.1428
$ws.Range("M98").Value = -581.1428000000001
$ws.Range("H111").Value = 1472.5
$ws.Range("I111").Value = 1392.5714
$ws.Range("K111").Value = 4177.7142
$ws.Range("M111").Value = -1110.7142
$ws.Range("H112").Value = 2951.6667
$ws.Range("J112").Value = 3148
$ws.Range("L112").Value = 9444
$ws.Range("N112").Value = -11660
$ws.Range("H113").Value = 2949.375
$ws.Range("I113").Value = 2599.1667
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 2599.1667
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = 654.8332999999998
$ws.Range("N113").Value = -10508
$ws.Range("H122").Value = 1963.5
$ws.Range("I122").Value = 2079.1428
$ws.Range("K122").Value = 6237.428400000001
$ws.Range("M122").Value = -3787.428400000001
$ws.Range("H135").Value = 5728.769
$ws.Range("I135").Value = 696.2
$ws.Range("J135").Value = 8874.125
$ws.Range("K135").Value = 6265.8
$ws.Range("L135").Value = 79867.125
$ws.Range("M135").Value = -3730.8
$ws.Range("N135").Value = -84937.125

# ---- ARM sheet ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8219.379999999999
$ws.Range("I32").Value = 3583.923
$ws.Range("J32").Value = 17146.926
$ws.Range("K32").Value = 3583.923
$ws.Range("L32").Value = 17146.926
$ws.Range("M32").Value = -3296.923
$ws.Range("N32").Value = -17720.926
$ws.Range("H61").Value = 3453704.2
$ws.Range("I61").Value = 5698.5654
$ws.Range("K61").Value = 5698.5654
$ws.Range("M61").Value = -5486.5654
$ws.Range("H132").Value = 896203.2
$ws.Range("I132").Value = 1045195.9
$ws.Range("J132").Value = 2247.25
$ws.Range("K132").Value = 3135587.7
$ws.Range("L132").Value = 6741.75
$ws.Range("M132").Value = -3133057.7
$ws.Range("N132").Value = -11801.75
$ws.Range("H136").Value = 3453704.2
$ws.Range("I136").Value = 5698.5654
$ws.Range("K136").Value = 17095.6962
$ws.Range("M136").Value = -14545.6962

# ---- BSM sheet ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2903.5173
$ws.Range("I105").Value = 2761.4773
$ws.Range("J105").Value = 3349.9285
$ws.Range("K105").Value = 2761.4773
$ws.Range("L105").Value = 3349.9285
$ws.Range("M105").Value = -1014.4773
$ws.Range("N105").Value = -6843.9285

# ---- CRP sheet ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 79350.08
$ws.Range("I16").Value = 2400.125
$ws.Range("K16").Value = 2400.125
$ws.Range("M16").Value = -2113.125
$ws.Range("H22").Value = 1144.2632
$ws.Range("I22").Value = 231.54546
$ws.Range("K22").Value = 231.54546
$ws.Range("M22").Value = 118.45454
$ws.Range("H86").Value = 11587.857
$ws.Range("I86").Value = 14387.556
$ws.Range("J86").Value = 6548.4
$ws.Range("K86").Value = 14387.556
$ws.Range("L86").Value = 6548.4
$ws.Range("M86").Value = -13264.556
$ws.Range("N86").Value = -8794.4
$ws.Range("H89").Value = 11587.857
$ws.Range("I89").Value = 14387.556
$ws.Range("J89").Value = 6548.4
$ws.Range("K89").Value = 71937.78
$ws.Range("L89").Value = 32742
$ws.Range("M89").Value = -66321.78
$ws.Range("N89").Value = -43974
$ws.Range("H105").Value = 10270.583
$ws.Range("J105").Value = 3562.6667
$ws.Range("L105").Value = 3562.6667
$ws.Range("N105").Value = -7056.6667
$ws.Range("H111").Value = 69107.2
$ws.Range("J111").Value = 68974.5
$ws.Range("L111").Value = 68974.5
$ws.Range("N111").Value = -77154.5
$ws.Range("H113").Value = 79350.08
$ws.Range("I113").Value = 2400.125
$ws.Range("K113").Value = 2400.125
$ws.Range("M113").Value = -230.125
$ws.Range("H116").Value = 77208
$ws.Range("J116").Value = 74742
$ws.Range("L116").Value = 74742
$ws.Range("N116").Value = -83920
$ws.Range("H132").Value = 4264.385
$ws.Range("I132").Value = 4203.25
$ws.Range("K132").Value = 12609.75
$ws.Range("M132").Value = -10079.75
$ws.Range("H134").Value = 1713.091
$ws.Range("I134").Value = 1771.5555
$ws.Range("K134").Value = 5314.666499999999
$ws.Range("M134").Value = -2779.666499999999

# ---- CUL sheet ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 434.375
$ws.Range("I8").Value = 434.375
$ws.Range("K8").Value = 1303.125
$ws.Range("M8").Value = -1164.125
$ws.Range("H40").Value = 71.75
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H137").Value = 5880.567
$ws.Range("I137").Value = 1909.4615
$ws.Range("J137").Value = 8917.294
$ws.Range("K137").Value = 5728.3845
$ws.Range("L137").Value = 26751.882
$ws.Range("M137").Value = -628.3845000000001
$ws.Range("N137").Value = -36951.882

# ---- GSM sheet ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 71.03333000000001
$ws.Range("I2").Value = 48.958332
$ws.Range("J2").Value = 159.33333
$ws.Range("K2").Value = 48.958332
$ws.Range("L2").Value = 159.33333
$ws.Range("M2").Value = 64.041668
$ws.Range("N2").Value = -385.33333
$ws.Range("H80").Value = 6546.3076
$ws.Range("I80").Value = 4531.5
$ws.Range("J80").Value = 9770
$ws.Range("K80").Value = 4531.5
$ws.Range("L80").Value = 9770
$ws.Range("M80").Value = -3533.5
$ws.Range("N80").Value = -11766
$ws.Range("H83").Value = 6546.3076
$ws.Range("I83").Value = 4531.5
$ws.Range("J83").Value = 9770
$ws.Range("K83").Value = 22657.5
$ws.Range("L83").Value = 48850
$ws.Range("M83").Value = -17665.5
$ws.Range("N83").Value = -58834
$ws.Range("H135").Value = 73309.766
$ws.Range("J135").Value = 73309.766
$ws.Range("L135").Value = 73309.766
$ws.Range("N135").Value = -83449.766
$ws.Range("H136").Value = 66566.09
$ws.Range("J136").Value = 66566.09
$ws.Range("L136").Value = 199698.27
$ws.Range("N136").Value = -204798.27

# ---- LTW sheet ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 12240
$ws.Range("I7").Value = 20276.555
$ws.Range("J7").Value = 3198.875
$ws.Range("K7").Value = 20276.555
$ws.Range("L7").Value = 3198.875
$ws.Range("M7").Value = -20164.555
$ws.Range("N7").Value = -3422.875
$ws.Range("H16").Value = 1287.7142
$ws.Range("J16").Value = 8175.5
$ws.Range("L16").Value = 8175.5
$ws.Range("N16").Value = -8515.5
$ws.Range("H46").Value = 3913.0667
$ws.Range("I46").Value = 558.3333
$ws.Range("J46").Value = 4751.75
$ws.Range("K46").Value = 558.3333
$ws.Range("L46").Value = 4751.75
$ws.Range("M46").Value = -370.3333
$ws.Range("N46").Value = -5127.75
$ws.Range("H55").Value = 1518
$ws.Range("I55").Value = 1453.2
$ws.Range("K55").Value = 1453.2
$ws.Range("M55").Value = -1280.2
$ws.Range("H74").Value = 50221.5
$ws.Range("I74").Value = 50221.5
$ws.Range("K74").Value = 50221.5
$ws.Range("M74").Value = -49223.5
$ws.Range("H77").Value = 50221.5
$ws.Range("I77").Value = 50221.5
$ws.Range("K77").Value = 150664.5
$ws.Range("M77").Value = -145672.5
$ws.Range("H93").Value = 1636.1538
$ws.Range("I93").Value = 852.3333
$ws.Range("K93").Value = 852.3333
$ws.Range("M93").Value = 395.6667
$ws.Range("H126").Value = 12240
$ws.Range("I126").Value = 20276.555
$ws.Range("J126").Value = 3198.875
$ws.Range("K126").Value = 60829.665
$ws.Range("L126").Value = 9596.625
$ws.Range("M126").Value = -58359.665
$ws.Range("N126").Value = -14536.625
$ws.Range("H128").Value = 40000
$ws.Range("J128").Value = 40000
$ws.Range("L128").Value = 40000
$ws.Range("N128").Value = -49960
$ws.Range("H138").Value = 99999
$ws.Range("J138").Value = 99999
$ws.Range("L138").Value = 99999
$ws.Range("N138").Value = -110279
$ws.Range("H139").Value = 89999
$ws.Range("J139").Value = 89999
$ws.Range("L139").Value = 89999
$ws.Range("N139").Value = -100279

# ---- WVR sheet ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 8776188
$ws.Range("I132").Value = 10420796
$ws.Range("K132").Value = 31262388
$ws.Range("M132").Value = -31259858
